$wb = $excel.ActiveWorkbook

# ALC!row15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1028.0605
$ws.Range("I15").Value = 1028.0605
$ws.Range("K15").Value = 3084.1815
$ws.Range("M15").Value = -2915.1815

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 9092943
$ws.Range("J19").Value = 11113175
$ws.Range("L19").Value = 11113175
$ws.Range("N19").Value = -11113525

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7862.75
$ws.Range("I86").Value = 2375.5
$ws.Range("K86").Value = 2375.5
$ws.Range("M86").Value = -1252.5

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7862.75
$ws.Range("I89").Value = 2375.5
$ws.Range("K89").Value = 11877.5
$ws.Range("M89").Value = -6261.5

# ALC!row94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 6438
$ws.Range("I94").Value = 6438
$ws.Range("K94").Value = 6438
$ws.Range("M94").Value = -5987

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2074.394
$ws.Range("I137").Value = 1917.5555
$ws.Range("J137").Value = 2780.1667
$ws.Range("K137").Value = 5752.666499999999
$ws.Range("L137").Value = 8340.500100000001
$ws.Range("M137").Value = -3202.666499999999
$ws.Range("N137").Value = -13440.5001

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8118.6816
$ws.Range("I138").Value = 7422
$ws.Range("J138").Value = 8273.5
$ws.Range("K138").Value = 22266
$ws.Range("L138").Value = 24820.5
$ws.Range("M138").Value = -17126
$ws.Range("N138").Value = -35100.5

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5214.6816
$ws.Range("I32").Value = 4510.619
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 4510.619
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -4223.619
$ws.Range("N32").Value = -20574

# ARM!row43
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 38655.8
$ws.Range("J43").Value = 38296.332
$ws.Range("L43").Value = 38296.332
$ws.Range("N43").Value = -38922.332

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3044.3845
$ws.Range("I74").Value = 1757.9
$ws.Range("K74").Value = 1757.9
$ws.Range("M74").Value = -883.9000000000001

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3044.3845
$ws.Range("I77").Value = 1757.9
$ws.Range("K77").Value = 8789.5
$ws.Range("M77").Value = -4421.5

# ARM!row92
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 54999.75
$ws.Range("J92").Value = 54999.75
$ws.Range("L92").Value = 54999.75
$ws.Range("N92").Value = -59991.75

# ARM!row112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 28000
$ws.Range("J112").Value = 28000
$ws.Range("L112").Value = 28000
$ws.Range("N112").Value = -30954

# ARM!row119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 79999.5
$ws.Range("J119").Value = 79999.5
$ws.Range("L119").Value = 79999.5
$ws.Range("N119").Value = -89675.5

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2503807.8
$ws.Range("I132").Value = 3778.0278
$ws.Range("J132").Value = 25004074
$ws.Range("K132").Value = 11334.0834
$ws.Range("L132").Value = 75012222
$ws.Range("M132").Value = -8804.0834
$ws.Range("N132").Value = -75017282

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27780886
$ws.Range("I31").Value = 35717170
$ws.Range("J31").Value = 3898.125
$ws.Range("K31").Value = 35717170
$ws.Range("L31").Value = 3898.125
$ws.Range("M31").Value = -35716875
$ws.Range("N31").Value = -4488.125

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 27780886
$ws.Range("I34").Value = 35717170
$ws.Range("J34").Value = 3898.125
$ws.Range("K34").Value = 35717170
$ws.Range("L34").Value = 3898.125
$ws.Range("M34").Value = -35716968
$ws.Range("N34").Value = -4302.125

# CRP!row51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 27246.166
$ws.Range("I51").Value = 27044.75
$ws.Range("J51").Value = 27649
$ws.Range("K51").Value = 27044.75
$ws.Range("L51").Value = 27649
$ws.Range("M51").Value = -26308.75
$ws.Range("N51").Value = -29121

# CRP!row60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

# CRP!row61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 27246.166
$ws.Range("I61").Value = 27044.75
$ws.Range("J61").Value = 27649
$ws.Range("K61").Value = 27044.75
$ws.Range("L61").Value = 27649
$ws.Range("M61").Value = -26696.75
$ws.Range("N61").Value = -28345

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2779.52
$ws.Range("I132").Value = 2651.652
$ws.Range("K132").Value = 7954.956
$ws.Range("M132").Value = -5424.956

# CUL!row12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1998.2142
$ws.Range("I12").Value = 81.40000000000001
$ws.Range("J12").Value = 3063.111
$ws.Range("K12").Value = 244.2
$ws.Range("L12").Value = 9189.332999999999
$ws.Range("M12").Value = -71.20000000000002
$ws.Range("N12").Value = -9535.332999999999

# CUL!row26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 4460.25
$ws.Range("I26").Value = 269.8
$ws.Range("J26").Value = 11444.333
$ws.Range("K26").Value = 809.4000000000001
$ws.Range("L26").Value = 34332.999
$ws.Range("M26").Value = -521.4000000000001
$ws.Range("N26").Value = -34908.999

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2188.8
$ws.Range("I80").Value = 1802.2
$ws.Range("K80").Value = 1802.2
$ws.Range("M80").Value = -804.2

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2188.8
$ws.Range("I83").Value = 1802.2
$ws.Range("K83").Value = 9011
$ws.Range("M83").Value = -4019

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4104.5884
$ws.Range("I102").Value = 4068.6
$ws.Range("K102").Value = 4068.6
$ws.Range("M102").Value = -2446.6

# GSM!row103
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 129949.5
$ws.Range("J103").Value = 129949.5
$ws.Range("L103").Value = 129949.5
$ws.Range("N103").Value = -132293.5

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2186
$ws.Range("I122").Value = 2267.111
$ws.Range("J122").Value = 2003.5
$ws.Range("K122").Value = 6801.333
$ws.Range("L122").Value = 6010.5
$ws.Range("M122").Value = -4351.333
$ws.Range("N122").Value = -10910.5

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1537.25

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4622.05
$ws.Range("J7").Value = 4273
$ws.Range("L7").Value = 4273
$ws.Range("N7").Value = -4497

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 930.19354
$ws.Range("I55").Value = 561.38464
$ws.Range("J55").Value = 1196.5555
$ws.Range("K55").Value = 561.38464
$ws.Range("L55").Value = 1196.5555
$ws.Range("M55").Value = -388.38464
$ws.Range("N55").Value = -1542.5555

# LTW!row63
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 99562.375
$ws.Range("J63").Value = 99562.375
$ws.Range("L63").Value = 99562.375
$ws.Range("N63").Value = -101060.375

# LTW!row66
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 99562.375
$ws.Range("J66").Value = 99562.375
$ws.Range("L66").Value = 298687.125
$ws.Range("N66").Value = -306175.125

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5912.143
$ws.Range("J82").Value = 6998.3335
$ws.Range("L82").Value = 6998.3335
$ws.Range("N82").Value = -7720.3335

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 5912.143
$ws.Range("J85").Value = 6998.3335
$ws.Range("L85").Value = 6998.3335
$ws.Range("N85").Value = -9494.333500000001

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4622.05
$ws.Range("J126").Value = 4273
$ws.Range("L126").Value = 12819
$ws.Range("N126").Value = -17759

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5722.5713
$ws.Range("I132").Value = 3333
$ws.Range("K132").Value = 9999
$ws.Range("M132").Value = -7469

# WVR!row15
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 37437.5
$ws.Range("I15").Value = 39900
$ws.Range("K15").Value = 39900
$ws.Range("M15").Value = -39612

# WVR!row96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5428.467
$ws.Range("I96").Value = 6741.2856
$ws.Range("J96").Value = 4279.75
$ws.Range("K96").Value = 6741.2856
$ws.Range("L96").Value = 4279.75
$ws.Range("M96").Value = -5368.2856
$ws.Range("N96").Value = -7025.75

# WVR!row100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1325.3125
$ws.Range("I100").Value = 1300.3334
$ws.Range("K100").Value = 2600.6668
$ws.Range("M100").Value = -2059.6668

# WVR!row119
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 89999
$ws.Range("J119").Value = 89999
$ws.Range("L119").Value = 89999
$ws.Range("N119").Value = -99675

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8259.223
$ws.Range("I126").Value = 8259.223
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 24777.669
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -22307.669
$ws.Range("N126").ClearContents()

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3335666.8
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 5002500
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 15007500
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -15012560
